$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Halves the target''s [Defense]{mechanic:defense} for damage calculation, which is similar to doubling the attack''s [power]{mechanic:power}."
$ws.Range("D17").Value = "Determines damage using each participant''s base [Attack]{mechanic:attack} and the opponent''s base [Defense]{mechanic:defense} in place of calculated stats."
$ws.Range("D18").Value = "[]{move:swift} can hit the user while in the air, and no other moves can.  Due to a bug, if the user fails to carry out the move''s second turn, e.g. due to [paralysis]{mechanic:paralysis}, it remains in the air until it switches out or uses Fly again successfully.  It can act as if it were on the ground in the meantime, but moves still miss it."
$ws.Range("D20").Value = "Does not raise the user''s [Special Defense]{mechanic:special-defense}."
$ws.Range("D23").Value = "[]{move:swift} can hit the user while underground, and no other moves can.  Due to a bug, if the user fails to carry out the move''s second turn, e.g. due to [paralysis]{mechanic:paralysis}, it remains underground until it switches out or uses Dig again successfully.  It can act as if it were aboveground in the meantime, but moves still miss it."
